$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set purpose (column E) and libraryPreparer (column B) for all data rows (2-43)
$ws.Range("E2:E43").Value = "fullRNASEQ"
$ws.Range("B2:B43").Value = "H.BROWN"

# Reflect the selection left behind by the editor in the saved view state
$ws.Range("B3:B43").Select()
